$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $id = $cell.Value()
    $cell.Value = $id + "(K)"
}

$ws.Range("A12").Select()
